$wb = $excel.ActiveWorkbook

$wsYearly = $wb.Worksheets.Item("Yearly")
$wsAllTime = $wb.Worksheets.Item("All Time")

# Update the 2017 (Yearly sheet) Jan figures
$wsYearly.Range("L3").Value = 50.57
$wsYearly.Range("M3").Value = 12.74
$wsYearly.Range("N3").Value = 5.26

# Update the All Time sheet corresponding year row (2017)
$wsAllTime.Range("G8").Value = 12.74
$wsAllTime.Range("H8").Value = 5.26

# Update selections to match final state
$wsYearly.Activate()
$wsYearly.Range("N4").Select()

$wsAllTime.Activate()
$excel.ActiveWindow.ScrollRow = 36
$wsAllTime.Range("K16").Select()

$wb.Save()
